$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the missile identifier used throughout the data rows (B:D) ---
# Shared strings: "Missile_HIGHWIND2_State_Update" -> "Missile_ANGERMAX2_State_Update"
#                 "MISSILE_HIGHWIND2_306.MISSILE_HIGHWIND2_306" -> "MISSILE_ANGERMAX2_468.MISSILE_ANGERMAX2_468"
#                 "MISSILE_HIGHWIND2" -> "MISSILE_ANGERMAX2"
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 2).Value = "Missile_ANGERMAX2_State_Update"
    $ws.Cells.Item($r, 3).Value = "MISSILE_ANGERMAX2_468.MISSILE_ANGERMAX2_468"
    $ws.Cells.Item($r, 4).Value = "MISSILE_ANGERMAX2"
}

# --- Updated simulation output numbers (columns F..K) for rows 2..54 ---
# Each inner array is: row, F, G, H, I, J, K ($null => leave the existing value untouched)
$data = @(
    @(2, $null, $null, $null, -1573.587032882787, 2032.167819573139, $null),
    @(3, $null, $null, $null, -1455.179115140841, 1974.576415652925, 316.2800275016236),
    @(4, $null, $null, $null, -1434.1253606197, 1965.595796181476, 583.1639455503765),
    @(5, $null, $null, $null, -1417.402595217621, 1851.395202956704, 848.6903663704209),
    @(6, $null, $null, $null, -1428.054747359177, 1844.057257215934, 1130.040382356711),
    @(7, $null, $null, $null, -1407.631118172161, 1786.859289519137, 1322.405336220726),
    @(8, $null, $null, $null, -1284.854029251887, 1673.859871536685, 1525.360931848673),
    @(9, $null, -104.9032508153096, $null, -1319.334413323036, 1706.634834989053, 1875.692663536948),
    @(10, 213.4861827150019, -83.8235347744374, 878.6023749890238, -1280.302969319163, 1583.506362198231, 1942.979460160965),
    @(11, 170.129049217227, -64.89130889164315, 1053.516959354966, -1149.707600775173, 1558.039237948055, 2157.033298989687),
    @(12, 135.5025085979326, -50.81158378858586, 1138.818869393986, -1132.161097084835, 1551.861012288439, 2296.485779119064),
    @(13, 125.9777894241739, -35.52017499769057, 1201.26963436906, -1133.120684844117, 1491.83812608717, 2588.708153379439),
    @(14, 108.8137119899331, -17.93265580991357, 1373.046654434757, -1048.357573851607, 1374.603933296502, 2542.237801879255),
    @(15, 97.16514177430021, -0.9272430183320194, 1433.209937748604, -1079.864639739991, 1412.118024637207, 2692.415232840783),
    @(16, 95.0376479819716, 15.04665370881238, 1473.757472820108, -1029.855452813576, 1358.765190820968, 2752.217806083632),
    @(17, 83.86022768392645, 32.43582364564039, 1517.484774530964, -972.4725874756493, 1296.700922999422, 2830.80735934923),
    @(18, 84.0914106364835, 50.60315086276741, 1517.286868442979, -930.3109081275181, 1265.287044781204, 2967.366618984558),
    @(19, 77.14993120963955, 68.36065919075436, 1609.370432021196, -917.9196239894678, 1146.81149391931, 3080.519679385576),
    @(20, 72.56724459481697, 80.48910240999909, 1624.304792001691, -835.5675167022373, 1170.684548364667, 3305.537146197419),
    @(21, 73.4220561295961, 102.2675630366356, 1519.438986011786, -826.2366436796773, 1037.549698869207, 3162.39779882563),
    @(22, 67.40996575937851, 109.4615750148627, 1562.499800235611, -739.4842570105152, 998.6149675804322, 3326.341842740402),
    @(23, 64.83081221207055, 124.7042453560242, 1665.042612561411, -728.0077785379316, 975.0441390031242, 3159.744436085696),
    @(24, 63.1802633891307, 150.1183885284138, 1659.004501227952, -668.9295348875107, 968.8021631323185, 3262.298105924876),
    @(25, 61.63080937766918, 163.5325562138718, 1608.356322625348, -604.8531176186052, 884.4884851397809, 3013.043544235238),
    @(26, 62.58686461344141, 178.8615675588663, 1639.750305612225, -564.7562781827546, 852.5463912914187, 3013.862997564465),
    @(27, 56.26666101855694, 191.9831184962643, 1778.909223390443, -469.2806652837993, 750.7508224323747, 3013.622284439266),
    @(28, 55.24259954684857, 205.823114049112, 1763.919186510255, -431.5266897928079, 756.6536976333118, 3022.379307873263),
    @(29, 55.7312258905984, 227.1513508602628, 1717.069377541605, -384.9474775797055, 697.4519118298657, 2971.474591042006),
    @(30, 54.2009096789393, 247.2022455247485, 1756.684406085497, -326.0826901213546, 643.3520015533271, 2865.962401643062),
    @(31, 53.33409439551652, 258.4038732918999, 1843.475823401016, -259.1936021000882, 582.1825361082961, 2750.254298656001),
    @(32, 51.14893594462347, 278.1865803335052, 1846.345310508557, -211.5455727023748, 538.5599070818739, 2523.042089744286),
    @(33, 47.86013926859631, 296.2702573734949, 1735.513290767555, -143.0867045217527, 474.4350986069132, 2316.70848125954),
    @(34, 48.08447459670824, 320.7417843608654, 1874.890501178554, -81.55379439407463, 445.8248920925536, 2127.869655760593),
    @(35, 50.11461017445732, 339.6937687209411, 1827.072596905194, -18.8021453904238, 388.6882307881807, 2109.321354599191),
    @(36, 48.28831298898171, 329.1138020145842, 1879.760248352524, 47.81067205418263, 355.7863393530487, 1885.795265415969),
    @(37, 47.79981952004732, 379.0314156394633, 1924.86257015555, 109.9941218001605, 284.4640032178323, 1584.231826597959),
    @(38, 47.6639504378576, 391.4796257987088, 1823.285224249025, 177.7184092596517, 238.251611322319, 1401.195069036983),
    @(39, 43.33049225929747, 392.8997930222339, 1951.376746322633, 258.4524330146502, 198.6479027168259, 1130.105288927545),
    @(40, 42.13801051890773, 423.4515587681851, 1924.205420819199, 324.4187085382803, 146.3572481880372, 923.7650330577743),
    @(41, 41.88782168883781, 437.073735378574, 1943.937846658554, 389.979689438239, 104.6519145021644, 641.2664451964006),
    @(42, 41.54909578175238, 448.7645098379222, 1857.696679748807, 460.5733813255522, 52.88020524242419, 350.3264467730222),
    @(43, 41.89318297376273, 456.7422771869212, 1817.072769628431, 542.6650117726352, 5.361397967670222, 34.25526678138262),
    @(44, 41.17604920802057, 479.5489392116456, 1998.206049079525, 612.0289528457338, -42.84872516041617, -285.7778212118015),
    @(45, 42.70933018626154, 482.1518819180341, 2010.659565174453, 726.0760812229324, -96.36146026552343, -600.1224731447538),
    @(46, 40.55609759288938, 489.2042498473633, 1986.903597143377, 797.164044024795, -143.0097966233286, -928.2759554518835),
    @(47, 40.95919789406019, 553.7571020219625, 1945.752895090712, 854.4317315580778, -182.1968072400078, -1364.526466426827),
    @(48, 37.52636691283905, 536.0614554507888, 1926.139173043984, 965.4963202540315, -239.5173826112327, -1619.802140270479),
    @(49, 38.66060522638427, 550.2982288724132, 1999.94183371625, 1074.276413224097, -292.0187023387205, -2033.572535229269),
    @(50, 38.98551191513928, 594.2045531338352, 1915.621131728892, 1198.05913088838, -329.7875744744806, -2626.990404955229),
    @(51, 37.56319192394465, 602.2286129535896, 1926.851098652049, 1268.506647824257, -392.962752026836, -2806.258746610545),
    @(52, 38.44200368650287, 615.0782367611706, 1978.751076569446, 1369.87656717986, -430.4182254480932, -3312.146965076264),
    @(53, 35.533018767313, 614.0807706851593, 1907.568128332946, 1465.528399736654, -505.0886882008357, -3916.622400610061),
    @(54, 36.95428884438494, 672.0749758684442, 2074.299614077109, 1543.578363771924, -516.7258378169637, -4284.99317888731)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 1; $i -le 6; $i++) {
        $val = $entry[$i]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $i + 5).Value = $val
        }
    }
}

Write-Host "edit applied"
